# Q factor data run for sg_rr_68_025 2023-12-11 17-15-27.csv
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$csvName   = "sg_rr_68_025 2023-12-11 17-15-27.csv"
$distanceNote = "(approx_fsr/2)/wavelength step size"

# ---- Row 83: first attempt at this run (prominence 0.005) ----
$ws.Cells.Item(83, 1).Value  = $csvName              # A83
$ws.Cells.Item(83, 2).Value  = 0.01                  # B83
$ws.Cells.Item(83, 3).Value  = 1000                  # C83
$ws.Cells.Item(83, 4).Value  = 5001                  # D83
$ws.Cells.Item(83, 5).Value  = 1530                  # E83
$ws.Cells.Item(83, 6).Value  = 1570                  # F83
$ws.Cells.Item(83, 7).Value  = 0.005                 # G83 - prominence/mW
$ws.Cells.Item(83, 8).Value  = $distanceNote         # H83 - distance
$ws.Cells.Item(83, 9).Value  = 1.7                   # I83 - approx_fsr/nm
$ws.Cells.Item(83, 21).Value = "found one peak in what looked like noise at end, so increased prominence slightly"  # U83

# ---- Row 84: re-run with slightly higher prominence (0.006) ----
$ws.Cells.Item(84, 1).Value  = $csvName              # A84
$ws.Cells.Item(84, 2).Value  = 0.01                  # B84
$ws.Cells.Item(84, 3).Value  = 1000                  # C84
$ws.Cells.Item(84, 4).Value  = 5001                  # D84
$ws.Cells.Item(84, 5).Value  = 1530                  # E84
$ws.Cells.Item(84, 6).Value  = 1570                  # F84
$ws.Cells.Item(84, 7).Value  = 0.006                 # G84 - prominence/mW
$ws.Cells.Item(84, 8).Value  = $distanceNote         # H84 - distance
$ws.Cells.Item(84, 9).Value  = 1.7                   # I84 - approx_fsr/nm
$ws.Cells.Item(84, 10).Value = 1.44807692307691      # J84 - fsr_mean/nm
$ws.Cells.Item(84, 11).Value = 0.00622696940163916   # K84 - fsr_std error/nm
$ws.Cells.Item(84, 12).Value = "yes, although  maybe misses peak at start which may be cut off in range"  # L84
$ws.Cells.Item(84, 13).Value = 0.145510029904756     # M84 - mean FWHM/nm
$ws.Cells.Item(84, 14).Value = 0.00407950032213657   # N84 - FWHM error/nm
$ws.Cells.Item(84, 15).Value = 10813.9579876443      # O84 - Q
$ws.Cells.Item(84, 16).Value = 225.456549156343      # P84 - Q error
$ws.Cells.Item(84, 17).Value = 273486698.550922      # Q84 - Q^3/R^2
$ws.Cells.Item(84, 18).Value = 17124396.8345199      # R84 - Q^3/R^2 error
$ws.Cells.Item(84, 19).Value = 68                    # S84 - radius/micrometres
$ws.Cells.Item(84, 20).Value = 0.1                   # T84 - radius error/micrometres

# ---- View state: scroll/zoom/selection to show the newly-added rows ----
$win = $excel.ActiveWindow
$win.Zoom = 105
$win.ScrollRow = 54
$win.ScrollColumn = 1
$ws.Range("A58").Select()
